# Generate Report for Handback
# Update the timestamp text values that record when handoff / handback
# xliff files were generated, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-25 07:04:59"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-25 07:04:53"
$wsZhCn.Range("K2").Value = "2016-08-25 07:05:19"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-25 07:04:59"
$wsDeDe.Range("K2").Value = "2016-08-25 07:05:26"
